$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition) - worksheet index 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 54
$ws1.Cells.Item(3, 6).Value = 11538
$ws1.Cells.Item(5, 6).Value = 326
$ws1.Cells.Item(7, 6).Value = 11488
$ws1.Cells.Item(8, 6).Value = 473
$ws1.Cells.Item(9, 6).Value = 1162
$ws1.Cells.Item(10, 6).Value = 81
$ws1.Cells.Item(11, 6).Value = 1755
$ws1.Cells.Item(12, 6).Value = 5700
$ws1.Cells.Item(13, 6).Value = 110
$ws1.Cells.Item(14, 6).Value = 3494
$ws1.Cells.Item(16, 6).Value = 14

# Sheet 4: 全部类型 (All types) - worksheet index 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 54
$ws4.Cells.Item(4, 6).Value = 3
$ws4.Cells.Item(5, 6).Value = 11538
$ws4.Cells.Item(7, 6).Value = 326
$ws4.Cells.Item(9, 6).Value = 11488
$ws4.Cells.Item(10, 6).Value = 473
$ws4.Cells.Item(11, 6).Value = 1162
$ws4.Cells.Item(12, 6).Value = 81
$ws4.Cells.Item(13, 6).Value = 1755
$ws4.Cells.Item(14, 6).Value = 3
$ws4.Cells.Item(15, 6).Value = 5700
$ws4.Cells.Item(16, 6).Value = 110
$ws4.Cells.Item(17, 6).Value = 3494
$ws4.Cells.Item(19, 6).Value = 14
